$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.985.07"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "1.639.83"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.35%  "
$ws.Range("D5").Value = "'215.04"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("D6").Value = "'0.5113"
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("D8").Value = "'0.2582"
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("D9").Value = "'0.06360"
$ws.Range("E9").Value = "  -0.83%  "
$ws.Range("D10").Value = "'19.80"
$ws.Range("E10").Value = "  +0.58%  "
$ws.Range("D11").Value = "'0.07764"
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").Value = "'4.281"
$ws.Range("E12").Value = "  -0.65%  "
$ws.Range("D13").Value = "1.629.59"
$ws.Range("E13").Value = "  -1.12%  "
$ws.Range("D14").Value = "'0.5465"
$ws.Range("E14").Value = "  -0.03%  "
$ws.Range("D15").Value = "0.0₅7750"
$ws.Range("E15").Value = "  -1.73%  "
$ws.Range("D16").Value = "'64.38"
$ws.Range("E16").Value = "  -0.70%  "
$ws.Range("D17").Value = "26.009.17"
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").Value = "'196.89"
$ws.Range("E19").Value = "  -0.79%  "
$ws.Range("D20").Value = "'4.429"
$ws.Range("E20").Value = "  -0.24%  "
$ws.Range("D21").Value = "'9.921"
$ws.Range("E21").Value = "  -1.07%  "
$ws.Range("D22").Value = "'6.090"
$ws.Range("E22").Value = "  +0.58%  "
$ws.Range("D23").Value = "'1.002"
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("D24").Value = "'1.917"
$ws.Range("E24").Value = "  +3.50%  "
$ws.Range("D25").Value = "'142.88"
$ws.Range("E25").Value = "  +2.01%  "
$ws.Range("D26").Value = "'0.1227"
$ws.Range("E26").Value = "  +7.07%  "
$ws.Range("D27").Value = "'6.853"
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("D28").Value = "'15.63"
$ws.Range("E28").Value = "  -0.77%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").Value = "'0.04852"
$ws.Range("E30").Value = "  -3.27%  "
$ws.Range("D31").Value = "'3.278"
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("D32").Value = "'3.214"
$ws.Range("E32").Value = "  +0.37%  "
$ws.Range("D33").Value = "'1.541"
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("E34").Value = "  +0.53%  "
$ws.Range("D35").Value = "'0.9152"
$ws.Range("E35").Value = "  +2.44%  "
$ws.Range("D36").Value = "'2.563"
$ws.Range("E36").Value = "  -0.64%  "
$ws.Range("D37").Value = "'0.5548"
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("D38").Value = "1.094.66"
$ws.Range("E38").Value = "  -3.38%  "
$ws.Range("D39").Value = "'0.01572"
$ws.Range("E39").Value = "  +0.54%  "
$ws.Range("E40").Value = "  -0.39%  "
$ws.Range("D41").Value = "'2.519"
$ws.Range("E41").Value = "  -1.31%  "
$ws.Range("D42").Value = "'5.574"
$ws.Range("E42").Value = "  -1.47%  "
$ws.Range("D43").Value = "'0.8050"
$ws.Range("E43").Value = "  -1.21%  "
$ws.Range("D44").Value = "'99.16"
$ws.Range("E44").Value = "  -0.60%  "
$ws.Range("E45").Value = "  -1.85%  "
$ws.Range("D46").Value = "1.777.19"
$ws.Range("E46").Value = "  -0.34%  "
$ws.Range("D47").Value = "'0.4533"
$ws.Range("E47").Value = "  +0.24%  "
$ws.Range("D48").Value = "'55.18"
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("D49").Value = "'1.005"
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").Value = "'0.05210"
$ws.Range("E50").Value = "  +2.34%  "
$ws.Range("D51").Value = "'7.520"
$ws.Range("E51").Value = "  +1.52%  "
